$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("NSDI")
$ws.Range("G2").Value = 'https://www.usenix.org/conference/nsdi23/presentation/cho-kun-woo'
$ws = $wb.Worksheets.Item("USENIX Security")
$ws.Range("G2").Value = 'https://www.usenix.org/conference/usenixsecurity22/presentation/shakevsky'
$ws.Range("G3").Value = 'https://www.usenix.org/conference/usenixsecurity22/presentation/gadotti'
$ws.Range("G4").Value = 'https://www.usenix.org/conference/usenixsecurity22/presentation/cerdeira'
$ws.Range("G5").Value = 'https://www.usenix.org/conference/usenixsecurity22/presentation/giner'
$ws.Range("G6").Value = 'https://www.usenix.org/conference/usenixsecurity22/presentation/cloosters'
$ws.Range("G7").Value = 'https://www.usenix.org/conference/usenixsecurity22/presentation/chen-yuan'
$ws.Range("G8").Value = 'https://www.usenix.org/conference/usenixsecurity22/presentation/kogler-minefield'
$ws.Range("G9").Value = 'https://www.usenix.org/conference/usenixsecurity23/presentation/mehmedagic'
$ws.Range("G10").Value = 'https://www.usenix.org/conference/usenixsecurity23/presentation/hilton'
$ws.Range("G11").Value = 'https://www.usenix.org/conference/usenixsecurity23/presentation/zhao-shixuan'
$ws.Range("G12").Value = 'https://www.usenix.org/conference/usenixsecurity23/presentation/constable'
$ws.Range("G13").Value = 'https://www.usenix.org/conference/usenixsecurity23/presentation/chen-hongbo'
$ws.Range("G14").Value = 'https://www.usenix.org/conference/usenixsecurity23/presentation/blechschmidt'
$ws.Range("G15").Value = 'https://www.usenix.org/conference/usenixsecurity23/presentation/zhang-haibin'
$ws.Range("G16").Value = 'https://www.usenix.org/conference/usenixsecurity23/presentation/caulfield'
$ws.Range("G17").Value = 'https://www.usenix.org/conference/usenixsecurity23/presentation/sass'
$ws.Range("G18").Value = 'https://www.usenix.org/conference/usenixsecurity24/presentation/wyss'
$ws.Range("G19").Value = 'https://www.usenix.org/conference/usenixsecurity24/presentation/schwarz'
$ws.Range("G20").Value = 'https://www.usenix.org/conference/usenixsecurity24/presentation/sridhara'
$ws.Range("G21").Value = 'https://www.usenix.org/conference/usenixsecurity24/presentation/schl%C3%BCter'
$ws.Range("G22").Value = 'https://www.usenix.org/conference/usenixsecurity24/presentation/busch-globalconfusion'
$ws = $wb.Worksheets.Item("NDSS")
$ws.Range("G2").Value = 'https://www.ndss-symposium.org/ndss-paper/mytee-own-the-trusted-execution-environment-on-embedded-devices/'
$ws.Range("G3").Value = 'https://www.ndss-symposium.org/ndss-paper/rr-a-fault-model-for-efficient-tee-replication/'
$ws.Range("G4").Value = 'https://www.ndss-symposium.org/ndss-paper/ldr-secure-and-efficient-linux-driver-runtime-for-embedded-tee-systems/'
$ws.Range("G5").Value = 'https://www.ndss-symposium.org/ndss-paper/overconfidence-is-a-dangerous-thing-mitigating-membership-inference-attacks-by-enforcing-less-confident-prediction/'
$ws.Range("G6").Value = 'https://www.ndss-symposium.org/ndss-paper/enclavefuzz-finding-vulnerabilities-in-sgx-applications/'
$ws.Range("G7").Value = 'https://www.ndss-symposium.org/ndss-paper/faults-in-our-bus-novel-bus-fault-attack-to-break-arm-trustzone/'
$ws.Range("G8").Value = 'https://www.ndss-symposium.org/ndss-paper/sense-enhancing-microarchitectural-awareness-for-tees-via-subscription-based-notification/'
$ws.Range("G9").Value = 'https://www.ndss-symposium.org/ndss-paper/tee-shirt-scalable-leakage-free-cache-hierarchies-for-tees/'
$ws = $wb.Worksheets.Item("OSDI")
$ws.Range("G2").Value = 'https://www.usenix.org/conference/osdi22/presentation/li'
$ws.Range("G3").Value = 'https://www.usenix.org/conference/osdi23/presentation/ahmad'
$ws.Range("G4").Value = 'https://www.usenix.org/conference/osdi23/presentation/angel'
$ws.Range("G5").Value = 'https://www.usenix.org/conference/osdi23/presentation/zhou-ziqiao'
$ws.Range("G6").Value = 'https://www.usenix.org/conference/osdi24/presentation/zhou'
$ws = $wb.Worksheets.Item("VLDB")
$ws.Range("G2").Value = 'https://ceur-ws.org/Vol-3462/DEco1.pdf'
$ws = $wb.Worksheets.Item("USENIX ATC")
$ws.Range("G2").Value = 'https://www.usenix.org/conference/atc22/presentation/sang'
$ws.Range("G3").Value = 'https://www.usenix.org/conference/atc23/presentation/li-dingji'
$ws.Range("G4").Value = 'https://www.usenix.org/conference/atc23/presentation/vaswani'
$ws.Range("G5").Value = 'https://www.usenix.org/conference/atc24/presentation/egorov'
$ws.Range("G6").Value = 'https://www.usenix.org/conference/atc24/presentation/chen-jiahao'
